# Update the "Idiomas" (Languages) worksheet with more descriptive
# video-download error messages.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlCenter = -4108

# ---------------------------------------------------------------------
# Row 24 stays the same content, only its height shrinks (45 -> 30)
# ---------------------------------------------------------------------
$ws.Rows.Item(24).RowHeight = 30

# ---------------------------------------------------------------------
# Row 25 is repurposed for the "age restricted" message (height 60 -> 30)
# ---------------------------------------------------------------------
$ws.Range("A25").Value() = "Este video tiene restricción de edad. `nNo puedes descargarlo."
$ws.Range("B25").Value() = "This video is age restricted. `nYou can't download it."
$ws.Range("A25:B25").HorizontalAlignment = $xlCenter
$ws.Range("A25:B25").WrapText = $True
$ws.Rows.Item(25).RowHeight = 30

# ---------------------------------------------------------------------
# Row 26 (new) - live stream still live message, same style as row 25
# ---------------------------------------------------------------------
$ws.Range("A26").Value() = "La transmisión en directo todavía esta vigente. `nVuelve a intentarlo cuando haya finalizado"
$ws.Range("B26").Value() = "The live stream is still live. `nPlease try again when it has finished"
$ws.Range("A26:B26").HorizontalAlignment = $xlCenter
$ws.Range("A26:B26").WrapText = $True
$ws.Rows.Item(26).RowHeight = 30

# ---------------------------------------------------------------------
# Rows 27-31 (new) - single line video status messages, style like row 23
# ---------------------------------------------------------------------
$ws.Range("A27").Value() = "Este video solo está disponible para miembros."
$ws.Range("B27").Value() = "This video is only available to members."

$ws.Range("A28").Value() = "Este video es privado."
$ws.Range("B28").Value() = "This video is private."

$ws.Range("A29").Value() = "Este video está bloqueado en tu región."
$ws.Range("B29").Value() = "This video is blocked in your region."

$ws.Range("A30").Value() = "Este video no está disponible."
$ws.Range("B30").Value() = "This video is unavailable."

$ws.Range("A31").Value() = " Ocurrió un error al descargar el video. Inténtalo nuevamente más tarde."
$ws.Range("B31").Value() = "There was an error downloading the video. Try again later."

$ws.Range("A27:B31").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# Row 32 (new) - blank row, only B32 carries formatting (no value)
# ---------------------------------------------------------------------
$ws.Range("B32").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# Column B is a touch wider now
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 56

# ---------------------------------------------------------------------
# Scroll / selection state
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("A25").Select() | Out-Null

Write-Host "Done"
